# Update "Ciudades" sheet in provincias_spain workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 22:46"

# Row 5 = Cataluña: update Casos totales (B5), Recuperados (D5), Muertes (E5)
$ws.Range("B5").Value = 4203
$ws.Range("D5").Value = 4078
$ws.Range("E5").Value = 122

# Row 29 = Tenerife: update Recuperados (D29), Muertes (E29)
$ws.Range("D29").Value = 184
$ws.Range("E29").Value = 4
